## Add a "last updated" date stamp to the About sheet (cell C1).
## The source workbook was re-saved from a newer Excel build; the only
## substantive content change is this new date value next to the title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Build a pure date (no time-of-day component) for 2021-04-21, which is
# serial day 44307 in the 1900 date system used by this workbook.
$dateValue = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$cell = $ws.Range("C1")
$cell.Value = $dateValue
$cell.NumberFormat = "m/d/yyyy"
